$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I18").Value = 5352.2383
$ws.Range("H18").Value = 5245.5835
$ws.Range("J18").Value = 4499
$ws.Range("K18").Value = 5352.2383
$ws.Range("M18").Value = -5068.2383
$ws.Range("L18").Value = 4499
$ws.Range("N18").Value = -5067
$ws.Range("J57").Value = 117666.336
$ws.Range("N57").Value = -353997.008
$ws.Range("L57").Value = 352999.008
$ws.Range("H57").Value = 117666.336
$ws.Range("I62").Value = 3333.5557
$ws.Range("H62").Value = 20349.55
$ws.Range("J62").Value = 34271.727
$ws.Range("K62").Value = 3333.5557
$ws.Range("M62").Value = -2709.5557
$ws.Range("L62").Value = 34271.727
$ws.Range("N62").Value = -35519.727
$ws.Range("I65").Value = 3333.5557
$ws.Range("J65").Value = 34271.727
$ws.Range("K65").Value = 16667.7785
$ws.Range("M65").Value = -13547.7785
$ws.Range("N65").Value = -177598.635
$ws.Range("L65").Value = 171358.635
$ws.Range("H65").Value = 20349.55
$ws.Range("J98").Value = 20682.572
$ws.Range("K98").Value = 28075.4
$ws.Range("M98").Value = -26577.4
$ws.Range("L98").Value = 20682.572
$ws.Range("N98").Value = -23678.572
$ws.Range("H98").Value = 26158.74
$ws.Range("I98").Value = 28075.4
$ws.Range("N106").Value = -11240.667
$ws.Range("L106").Value = 9978.666999999999
$ws.Range("H106").Value = 10293755
$ws.Range("J106").Value = 9978.666999999999
$ws.Range("N122").Value = -66947.716
$ws.Range("I122").Value = 28075.4
$ws.Range("H122").Value = 26158.74
$ws.Range("J122").Value = 20682.572
$ws.Range("K122").Value = 84226.20000000001
$ws.Range("M122").Value = -81776.20000000001
$ws.Range("L122").Value = 62047.716
$ws.Range("I132").Value = 5290.737
$ws.Range("H132").Value = 4986.72
$ws.Range("K132").Value = 15872.211
$ws.Range("M132").Value = -13342.211
$ws.Range("M137").Value = -34756.05
$ws.Range("J137").Value = 4946.5
$ws.Range("N137").Value = -19939.5
$ws.Range("L137").Value = 14839.5
$ws.Range("I137").Value = 12435.35
$ws.Range("H137").Value = 11754.546
$ws.Range("K137").Value = 37306.05
$ws.Range("H138").Value = 188800.56
$ws.Range("K138").Value = 1466361.18
$ws.Range("J138").Value = 4717.9316
$ws.Range("M138").Value = -1461221.18
$ws.Range("N138").Value = -24433.7948
$ws.Range("L138").Value = 14153.7948
$ws.Range("I138").Value = 488787.06

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M5").Value = -1688.8889
$ws.Range("N5").Value = -5204
$ws.Range("L5").Value = 4980
$ws.Range("H5").Value = 2779.077
$ws.Range("I5").Value = 1800.8889
$ws.Range("J5").Value = 4980
$ws.Range("K5").Value = 1800.8889
$ws.Range("J32").Value = 11996.5
$ws.Range("N32").Value = -12570.5
$ws.Range("L32").Value = 11996.5
$ws.Range("H32").Value = 2679.7942
$ws.Range("I45").Value = 3719.1538
$ws.Range("H45").Value = 4991.8945
$ws.Range("K45").Value = 3719.1538
$ws.Range("M45").Value = -3342.1538
$ws.Range("M61").Value = -5393.0303
$ws.Range("I61").Value = 5605.0303
$ws.Range("H61").Value = 5536.2104
$ws.Range("K61").Value = 5605.0303
$ws.Range("K74").Value = 1035
$ws.Range("M74").Value = -161
$ws.Range("I74").Value = 1035
$ws.Range("H74").Value = 4843.9287
$ws.Range("K77").Value = 5175
$ws.Range("M77").Value = -807
$ws.Range("I77").Value = 1035
$ws.Range("H77").Value = 4843.9287
$ws.Range("M110").Value = 353.8823
$ws.Range("I110").Value = 1691.1177
$ws.Range("H110").Value = 1997.5238
$ws.Range("K110").Value = 1691.1177
$ws.Range("I132").Value = 4579.052
$ws.Range("H132").Value = 4590.0806
$ws.Range("K132").Value = 13737.156
$ws.Range("M132").Value = -11207.156
$ws.Range("K136").Value = 16815.0909
$ws.Range("M136").Value = -14265.0909
$ws.Range("I136").Value = 5605.0303
$ws.Range("H136").Value = 5536.2104

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M4").Value = -1685.8889
$ws.Range("N4").Value = -5210
$ws.Range("J4").Value = 4980
$ws.Range("L4").Value = 4980
$ws.Range("I4").Value = 1800.8889
$ws.Range("H4").Value = 2779.077
$ws.Range("K4").Value = 1800.8889
$ws.Range("H11").Value = 7610.1113
$ws.Range("J11").Value = 5936
$ws.Range("K11").Value = 8949.4
$ws.Range("M11").Value = -8809.4
$ws.Range("L11").Value = 5936
$ws.Range("N11").Value = -6216
$ws.Range("I11").Value = 8949.4
$ws.Range("M20").Value = -799.5834
$ws.Range("N20").Value = -2761.4375
$ws.Range("L20").Value = 2267.4375
$ws.Range("H20").Value = 1744.2142
$ws.Range("I20").Value = 1046.5834
$ws.Range("J20").Value = 2267.4375
$ws.Range("K20").Value = 1046.5834
$ws.Range("J86").Value = 1990.3
$ws.Range("K86").Value = 4087.2917
$ws.Range("M86").Value = -2964.2917
$ws.Range("N86").Value = -4236.3
$ws.Range("L86").Value = 1990.3
$ws.Range("I86").Value = 4087.2917
$ws.Range("H86").Value = 3470.5293
$ws.Range("M89").Value = -14820.4585
$ws.Range("N89").Value = -21183.5
$ws.Range("L89").Value = 9951.5
$ws.Range("I89").Value = 4087.2917
$ws.Range("H89").Value = 3470.5293
$ws.Range("K89").Value = 20436.4585
$ws.Range("J89").Value = 1990.3
$ws.Range("K94").Value = 1093.4286
$ws.Range("M94").Value = -642.4286
$ws.Range("H94").Value = 2984.3333
$ws.Range("I94").Value = 1093.4286
$ws.Range("K99").Value = 16672.115
$ws.Range("M99").Value = -15174.115
$ws.Range("I99").Value = 16672.115
$ws.Range("H99").Value = 12830.77
$ws.Range("K105").Value = 2723.125
$ws.Range("M105").Value = -976.125
$ws.Range("I105").Value = 2723.125
$ws.Range("H105").Value = 3798.3928

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M12").Value = -231.33334
$ws.Range("I12").Value = 401.33334
$ws.Range("H12").Value = 643
$ws.Range("K12").Value = 401.33334
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("K99").Value = 23225912
$ws.Range("M99").Value = -23224414
$ws.Range("I99").Value = 23225912
$ws.Range("H99").Value = 23225912
$ws.Range("I126").Value = 23225912
$ws.Range("H126").Value = 23225912
$ws.Range("K126").Value = 69677736
$ws.Range("M126").Value = -69675266
$ws.Range("L139").Value = 99779.5
$ws.Range("H139").Value = 99779.5
$ws.Range("J139").Value = 99779.5
$ws.Range("N139").Value = -110059.5
$ws.Range("N52").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J113").Value = 1900.3334
$ws.Range("N113").Value = -10041.0002
$ws.Range("L113").Value = 5701.0002
$ws.Range("H113").Value = 1587.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M97").Value = -16495.357
$ws.Range("I97").Value = 16991.357
$ws.Range("H97").Value = 15110.294
$ws.Range("K97").Value = 16991.357
$ws.Range("L100").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("N106").Value = -1752524
$ws.Range("L106").Value = 1750000
$ws.Range("H106").Value = 1750000
$ws.Range("J106").Value = 1750000
$ws.Range("J113").Value = 3013
$ws.Range("K113").Value = 25359.8
$ws.Range("M113").Value = -23189.8
$ws.Range("L113").Value = 3013
$ws.Range("I113").Value = 25359.8
$ws.Range("H113").Value = 21635.334
$ws.Range("I126").Value = 24949.334
$ws.Range("H126").Value = 20780.52
$ws.Range("K126").Value = 74848.00199999999
$ws.Range("M126").Value = -72378.00199999999
$ws.Range("N113").Value = -7353
$ws.Range("N100").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I10").Value = 550
$ws.Range("H10").Value = 623
$ws.Range("K10").Value = 550
$ws.Range("M10").Value = -410
$ws.Range("N12").Value = -2473
$ws.Range("L12").Value = 2133
$ws.Range("I12").Value = 0
$ws.Range("H12").Value = 2133
$ws.Range("J12").Value = 2133
$ws.Range("K12").Value = 0
$ws.Range("L22").Value = 5508.3335
$ws.Range("N22").Value = -6098.3335
$ws.Range("I22").Value = 13031.6
$ws.Range("H22").Value = 8928
$ws.Range("J22").Value = 5508.3335
$ws.Range("K22").Value = 13031.6
$ws.Range("M22").Value = -12736.6
$ws.Range("J27").Value = 5508.3335
$ws.Range("K27").Value = 13031.6
$ws.Range("M27").Value = -12924.6
$ws.Range("L27").Value = 5508.3335
$ws.Range("N27").Value = -5722.3335
$ws.Range("I27").Value = 13031.6
$ws.Range("H27").Value = 8928
$ws.Range("M55").Value = -122.875
$ws.Range("N55").Value = -3278.1667
$ws.Range("J55").Value = 2932.1667
$ws.Range("L55").Value = 2932.1667
$ws.Range("I55").Value = 295.875
$ws.Range("H55").Value = 823.13336
$ws.Range("K55").Value = 295.875
$ws.Range("K93").Value = 11805.667
$ws.Range("M93").Value = -10557.667
$ws.Range("I93").Value = 11805.667
$ws.Range("H93").Value = 11189.77
$ws.Range("N122").Value = -64898.5
$ws.Range("I122").Value = 4905.8125
$ws.Range("H122").Value = 6582.8887
$ws.Range("J122").Value = 19999.5
$ws.Range("K122").Value = 14717.4375
$ws.Range("M122").Value = -12267.4375
$ws.Range("L122").Value = 59998.5
$ws.Range("I132").Value = 878088.4
$ws.Range("H132").Value = 575984.6
$ws.Range("K132").Value = 2634265.2
$ws.Range("M132").Value = -2631735.2
$ws.Range("M12").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I10").Value = 2305
$ws.Range("H10").Value = 3900.7144
$ws.Range("K10").Value = 2305
$ws.Range("M10").Value = -2136
$ws.Range("K13").Value = 1705.75
$ws.Range("M13").Value = -1565.75
$ws.Range("I13").Value = 1705.75
$ws.Range("H13").Value = 1125.5714
$ws.Range("I132").Value = 7210.364
$ws.Range("H132").Value = 6311.1143
$ws.Range("J132").Value = 4789.3076
$ws.Range("K132").Value = 21631.092
$ws.Range("M132").Value = -19101.092
$ws.Range("L132").Value = 14367.9228
$ws.Range("N132").Value = -19427.9228
